# Backlog.xlsx — mark radiator-related Hoja1 rows as "terminado" and note
# the unmatched row as "no se aplica" (per commit: identifying correlative
# radiator numbers that were superimposed / overlapping).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 13: status -> terminado; add note in column C
$ws.Range("B13").Value = "terminado"
$ws.Range("C13").Value = "no se aplica"

# Row 14: status -> terminado
$ws.Range("B14").Value = "terminado"

# Row 25: status -> terminado
$ws.Range("B25").Value = "terminado"

# Row 26: status -> terminado
$ws.Range("B26").Value = "terminado"

# Row 27: status -> terminado
$ws.Range("B27").Value = "terminado"

# Row 28: status -> terminado
$ws.Range("B28").Value = "terminado"

# Update the sheet view: zoom in and scroll/select as left by the editor.
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 115
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("B18").Select()
